$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q1" right after "2021-Q4" (and therefore
#    right before "总计").
# ---------------------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Item("2021-Q4")

$newSheet = $wb.Worksheets.Add($null, $wsQ4)
$newSheet.Name = "2022-Q1"

# Re-fetch sheet references by name now that the collection has changed, to
# avoid any stale-index issues after the insertion.
$wsQ1 = $wb.Worksheets.Item("2022-Q1")
$wsTotal = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# 2) Helper to stamp a cell with the same bold/boxed "header" look already
#    used throughout the workbook (header row + first data column). We copy
#    the cell *format only* from a cell that already has that exact look
#    (the "2021-Q4" sheet's B1 / A2 cells), which reproduces the existing
#    style instead of inventing a new one.
# ---------------------------------------------------------------------------
$headerStyleSource = $wsQ4.Range("B1")
$indexStyleSource = $wsQ4.Range("A2")

function Copy-FormatOnly($srcCell, $dstRange) {
    foreach ($cell in $dstRange.Cells) {
        $srcCell.Copy()
        $cell.PasteSpecial(-4122)
    }
}

# ---------------------------------------------------------------------------
# 3) Populate header row for "2022-Q1"
# ---------------------------------------------------------------------------
$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"
Copy-FormatOnly $headerStyleSource $wsQ1.Range("B1:H1")

# ---------------------------------------------------------------------------
# 4) Populate data rows for "2022-Q1"
# ---------------------------------------------------------------------------
$wsQ1.Range("A2").Value = 0
$wsQ1.Range("B2").Value = "'008851"
$wsQ1.Range("C2").Value = "景顺长城量化对冲策略三个月定期开放灵活配置混合"
$wsQ1.Range("D2").Value = "'5.05"
$wsQ1.Range("E2").Value = "'74.55"
$wsQ1.Range("F2").Value = "'1.38"
$wsQ1.Range("G2").Value = "'0.0697"
$wsQ1.Range("H2").Value = 10
Copy-FormatOnly $indexStyleSource $wsQ1.Range("A2")

$wsQ1.Range("A3").Value = 1
$wsQ1.Range("B3").Value = "'003704"
$wsQ1.Range("C3").Value = "光大保德信事件驱动灵活配置混合"
$wsQ1.Range("D3").Value = "'3.06"
$wsQ1.Range("E3").Value = "'23.55"
$wsQ1.Range("F3").Value = "'1.40"
$wsQ1.Range("G3").Value = "'0.0428"
$wsQ1.Range("H3").Value = 4
Copy-FormatOnly $indexStyleSource $wsQ1.Range("A3")

# ---------------------------------------------------------------------------
# 5) Update "总计" sheet: insert a new row for "2022-Q1" above the existing
#    "2021-Q4" row, shifting it down.
# ---------------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()
# Inserting a row copies the formatting of the row above (the bold header),
# so clear it from the data cells before writing the new values.
$wsTotal.Range("A2:D2").ClearFormats()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.11
Copy-FormatOnly $indexStyleSource $wsTotal.Range("A2")

$wsTotal.Range("A3").Value = 1

# ---------------------------------------------------------------------------
# 6) Restore the originally-active sheet/tab.
# ---------------------------------------------------------------------------
$wsQ4.Activate()
